$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$ppeSheet = $wb.Worksheets.Item("PPEIdtIL")

# --- Insert the new "Texas Notes" sheet between "About" and "PPEIdtIL" ---
$ws = $wb.Worksheets.Add()
$ws.Name = "Texas Notes"
$ws.Move($null, $aboutSheet)

# --- Text / value content (column A) ---
$ws.Range("A1").Value = "The study says:"

$ws.Range("A3").Value = "Table 2. Costs and benefits of improved appliance labeling"
$ws.Range("A4").Value = "Average % energy savings"
$ws.Range("A5").Value = 0.1
$ws.Range("A5").NumberFormat = "0%"

$ws.Range("A6").Value = "Assumes consumer selects products that average 10% energy savings. Difference in energy use from least to most efficient product varies from 10% to 50% depending on product category according to FTC published ranges. We assume most consumers motivated by the label select a product that falls 1-2 categories more efficient (e.g., shift from 1-star to 2- or 3-star product, or from 3-star to 4- or 5-star product, etc.). Source: ACEEE estimate based on labeling research."

$ws.Range("A11").Value = "% of market affected by policy"
$ws.Range("A12").Value = 0.2
$ws.Range("A12").NumberFormat = "0%"

$ws.Range("A13").Value = "Assumes 20% of consumers are motivated by the improved label and use it to select more efficient products. Source: ACEEE estimate based on prior labeling research"

$ws.Range("A16").Value = "Question for EI"
$ws.Range("A17").Value = "does this mean that the numbers in this spreadsheet should be 2% (10% * 20%)?"
$ws.Range("A18").Value = "right now they are 10%."

$ws.Range("A20").Value = "*Anwer from EI: yes, this was an error."

# --- Formatting: build up alignment styles row-group by row-group ---

# A11: left/top, no wrap
$ws.Range("A11").HorizontalAlignment = -4131
$ws.Range("A11").VerticalAlignment = -4160

# B11:M11: left/top + wrap
$ws.Range("B11:M11").HorizontalAlignment = -4131
$ws.Range("B11:M11").VerticalAlignment = -4160
$ws.Range("B11:M11").WrapText = $true

# A6:M10 and A13:M14: left/top + wrap (same combo as B11:M11)
$ws.Range("A6:M10").HorizontalAlignment = -4131
$ws.Range("A6:M10").VerticalAlignment = -4160
$ws.Range("A6:M10").WrapText = $true

$ws.Range("A13:M14").HorizontalAlignment = -4131
$ws.Range("A13:M14").VerticalAlignment = -4160
$ws.Range("A13:M14").WrapText = $true

# A15:M15: top + wrap (no horizontal)
$ws.Range("A15:M15").VerticalAlignment = -4160
$ws.Range("A15:M15").WrapText = $true

# A16:M17: top only
$ws.Range("A16:M17").VerticalAlignment = -4160

# A18:M23: touch wrap on/off to force an applyAlignment flag with defaults
$ws.Range("A18:M23").WrapText = $true
$ws.Range("A18:M23").WrapText = $false

# --- Merges (after formatting) ---
$ws.Range("A6:M10").Merge()
$ws.Range("A13:M14").Merge()

# --- Row height for row 13 (explicit 15pt) ---
$ws.Rows.Item(13).RowHeight = 15

Write-Output "done"
